# Weekly update: insert a new "Acelga" price record for "Vega Modelo de
# Temuco" at the top of the data block (row 294), pushing the existing
# rows 294-363 down to 295-364.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 294 - this shifts rows 294:363 down to 295:364
# (values, formatting and styles all move with the cells, matching the
# original rows exactly).
$ws.Rows.Item(294).Insert()

# Populate the freshly inserted row 294 with the new weekly data point.
$ws.Cells.Item(294, 1).Value  = 10                                  # A Mercado ID
$ws.Cells.Item(294, 2).Value  = "Vega Modelo de Temuco"              # B Mercado
$ws.Cells.Item(294, 3).Value  = "La Araucanía"                      # C Región
$ws.Cells.Item(294, 4).Value  = 44785                                # D Fecha
$ws.Cells.Item(294, 5).Value  = 9                                    # E Codreg
$ws.Cells.Item(294, 6).Value  = 100112009                            # F Categoría ID
$ws.Cells.Item(294, 7).Value  = "Acelga"                             # G Categoría
$ws.Cells.Item(294, 8).Value  = "Sin especificar"                    # H Variedad
$ws.Cells.Item(294, 9).Value  = "Primera"                            # I Calidad
$ws.Cells.Item(294, 10).Value = 30                                   # J Volumen
$ws.Cells.Item(294, 11).Value = 10000                                # K Precio mínimo
$ws.Cells.Item(294, 12).Value = 10000                                # L Precio máximo
$ws.Cells.Item(294, 13).Value = 10000                                # M Precio promedio ponderado
$ws.Cells.Item(294, 14).Value = "$/docena de atados (12 kilos)"      # N Unidad de comercialización
$ws.Cells.Item(294, 15).Value = "Provincia de Cautín"                # O Origen
$ws.Cells.Item(294, 16).Value = 833                                  # P Precio $/Kg
$ws.Cells.Item(294, 17).Value = 12                                   # Q Kg o Unidades
$ws.Cells.Item(294, 18).Value = "Hortaliza"                          # R Clasificación
